$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill season record values (Wins=71, Losses=91, Ties=0) for every player row
$ws.Range("AD2:AD54").Value = 71
$ws.Range("AE2:AE54").Value = 91
$ws.Range("AF2:AF54").Value = 0
